$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in row 2 to match the Birch+K_means clustering result
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 3

# Remove row 3 entirely (data no longer needed after clustering update)
$ws.Rows.Item(3).Delete()
